$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Row 25
$ws.Range("A25").Value = "DD OG AD 10(11) forsat"
$ws.Range("B25").Value = "Requirement Specifier"
$ws.Range("C25").Value = 43895
$ws.Range("D25").Value = 0.35416666666666669
$ws.Range("E25").Value = 0.41666666666666669

# Row 26
$ws.Range("A26").Value = "ATD 10"
$ws.Range("B26").Value = "Test Desinger"
$ws.Range("C26").Value = 43895
$ws.Range("D26").Value = 0.41666666666666669
$ws.Range("E26").Value = 0.47916666666666669

# Row 27
$ws.Range("A27").Value = "DCD0803"
$ws.Range("B27").Value = "Designer"
$ws.Range("C27").Value = 43895
$ws.Range("D27").Value = 0.52083333333333337
$ws.Range("E27").Value = 0.58333333333333337

# Row 28
$ws.Range("A28").Value = "SD 0803"
$ws.Range("B28").Value = "Designer"
$ws.Range("C28").Value = 43895
$ws.Range("D28").Value = 0.58333333333333337
$ws.Range("E28").Value = 0.65277777777777779

# View state update
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 14
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A29").Select()
